$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 14's formatting (incl. the date style on column A) into row 15
$ws.Range("A14:M14").Copy($ws.Range("A15:M15"))

# Fill in the actual values for the new row
$ws.Range("A15").Value = 43803
$ws.Range("B15").Value = 0
$ws.Range("C15").Value = 10
$ws.Range("D15").Value = 5
$ws.Range("E15").Value = 4
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 6
$ws.Range("H15").Value = 25
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 3

$ws.Range("H15").Select()
